{"js": "// Update the date line and all 25 division-problem cells in the table.\n// Each (oldText -> newText) pair is unique within the document, so we can\n// safely locate-and-replace each one independently via Body.search().\nconst pairs = [\n  [\"2024-02-19 Monday\", \"2024-02-20 Tuesday\"],\n  [\"320\u00f76=53, 2\", \"441\u00f74=110, 1\"],\n  [\"417\u00f79=46, 3\", \"958\u00f79=106, 4\"],\n  [\"723\u00f78=90, 3\", \"239\u00f76=39, 5\"],\n  [\"310\u00f73=103, 1\", \"475\u00f75=95, 0\"],\n  [\"947\u00f74=236, 3\", \"302\u00f78=37, 6\"],\n  [\"273\u00f74=68, 1\", \"622\u00f79=69, 1\"],\n  [\"885\u00f78=110, 5\", \"319\u00f78=39, 7\"],\n  [\"127\u00f74=31, 3\", \"605\u00f72=302, 1\"],\n  [\"874\u00f78=109, 2\", \"373\u00f77=53, 2\"],\n  [\"666\u00f75=133, 1\", \"637\u00f77=91, 0\"],\n  [\"419\u00f78=52, 3\", \"709\u00f78=88, 5\"],\n  [\"732\u00f78=91, 4\", \"930\u00f76=155, 0\"],\n  [\"383\u00f74=95, 3\", \"783\u00f76=130, 3\"],\n  [\"762\u00f74=190, 2\", \"685\u00f78=85, 5\"],\n  [\"479\u00f75=95, 4\", \"918\u00f73=306, 0\"],\n  [\"789\u00f78=98, 5\", \"689\u00f77=98, 3\"],\n  [\"316\u00f74=79, 0\", \"371\u00f73=123, 2\"],\n  [\"475\u00f76=79, 1\", \"574\u00f78=71, 6\"],\n  [\"929\u00f78=116, 1\", \"159\u00f73=53, 0\"],\n  [\"562\u00f72=281, 0\", \"931\u00f76=155, 1\"],\n  [\"270\u00f76=45, 0\", \"812\u00f77=116, 0\"],\n  [\"480\u00f78=60, 0\", \"144\u00f77=20, 4\"],\n  [\"950\u00f77=135, 5\", \"238\u00f72=119, 0\"],\n  [\"229\u00f72=114, 1\", \"523\u00f75=104, 3\"],\n  [\"598\u00f75=119, 3\", \"170\u00f78=21, 2\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  // Replace just the first (and expected only) match, preserving its run formatting.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and all 25 division-problem cells in the table.\n# Each (oldText -> newText) pair is unique within the document, so Find/Replace\n# One at a time (wdReplaceOne) safely targets exactly the intended run.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-02-19 Monday\", \"2024-02-20 Tuesday\"),\n    @(\"320\u00f76=53, 2\", \"441\u00f74=110, 1\"),\n    @(\"417\u00f79=46, 3\", \"958\u00f79=106, 4\"),\n    @(\"723\u00f78=90, 3\", \"239\u00f76=39, 5\"),\n    @(\"310\u00f73=103, 1\", \"475\u00f75=95, 0\"),\n    @(\"947\u00f74=236, 3\", \"302\u00f78=37, 6\"),\n    @(\"273\u00f74=68, 1\", \"622\u00f79=69, 1\"),\n    @(\"885\u00f78=110, 5\", \"319\u00f78=39, 7\"),\n    @(\"127\u00f74=31, 3\", \"605\u00f72=302, 1\"),\n    @(\"874\u00f78=109, 2\", \"373\u00f77=53, 2\"),\n    @(\"666\u00f75=133, 1\", \"637\u00f77=91, 0\"),\n    @(\"419\u00f78=52, 3\", \"709\u00f78=88, 5\"),\n    @(\"732\u00f78=91, 4\", \"930\u00f76=155, 0\"),\n    @(\"383\u00f74=95, 3\", \"783\u00f76=130, 3\"),\n    @(\"762\u00f74=190, 2\", \"685\u00f78=85, 5\"),\n    @(\"479\u00f75=95, 4\", \"918\u00f73=306, 0\"),\n    @(\"789\u00f78=98, 5\", \"689\u00f77=98, 3\"),\n    @(\"316\u00f74=79, 0\", \"371\u00f73=123, 2\"),\n    @(\"475\u00f76=79, 1\", \"574\u00f78=71, 6\"),\n    @(\"929\u00f78=116, 1\", \"159\u00f73=53, 0\"),\n    @(\"562\u00f72=281, 0\", \"931\u00f76=155, 1\"),\n    @(\"270\u00f76=45, 0\", \"812\u00f77=116, 0\"),\n    @(\"480\u00f78=60, 0\", \"144\u00f77=20, 4\"),\n    @(\"950\u00f77=135, 5\", \"238\u00f72=119, 0\"),\n    @(\"229\u00f72=114, 1\", \"523\u00f75=104, 3\"),\n    @(\"598\u00f75=119, 3\", \"170\u00f78=21, 2\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # MatchCase:true, MatchWholeWord:false, MatchWildcards:false, MatchSoundsLike:false,\n    # MatchAllWordForms:false, Forward:true, Wrap:wdFindContinue(1), Format:false,\n    # Replace:wdReplaceOne(2)\n    $result = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $result) {\n        throw \"Could not find text: $oldText\"\n    }\n}\n\n"}
